$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (player, position, team) in the desired final order
$data = @(
    @("Chris Paul", "PG", "San Antonio Spurs"),
    @("Stephon Castle", "PG,SG", "San Antonio Spurs"),
    @("Russell Westbrook", "PG,SG", "Denver Nuggets"),
    @("Cole Anthony", "PG", "Orlando Magic"),
    @("Jerami Grant", "SF,PF", "Portland Trail Blazers"),
    @("Kyle Filipowski", "PF,C", "Utah Jazz"),
    @("Paolo Banchero", "SF,PF", "Orlando Magic"),
    @("Jaylen Brown", "SG,SF", "Boston Celtics"),
    @("Pascal Siakam", "SF,PF,C", "Indiana Pacers"),
    @("Chet Holmgren", "PF,C", "Oklahoma City Thunder"),
    @("Nikola Jokic", "C", "Denver Nuggets"),
    @("Jalen Green", "PG,SG", "Houston Rockets"),
    @("Payton Pritchard", "PG,SG", "Boston Celtics"),
    @("Deni Avdija", "SF,PF", "Portland Trail Blazers"),
    @("Rudy Gobert", "C", "Minnesota Timberwolves"),
    @("Jakob Poeltl", "C", "Toronto Raptors")
)

# Clear out the old table body (previously rows 2-18) before writing the new, shorter table
$ws.Range("A2:C18").Clear()

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
